{"js": "// Add a new paragraph style \"byline\" based on the existing \"Author\"\n// style, matching the commit's styles.xml / stylesWithEffects.xml\n// addition:\n//\n//   <w:style w:type=\"paragraph\" w:customStyle=\"1\" w:styleId=\"byline\">\n//     <w:name w:val=\"byline\"/>\n//     <w:basedOn w:val=\"Author\"/>\n//     <w:qFormat/>\n//     <w:rsid w:val=\"00076B7A\"/>\n//   </w:style>\n\n// Create the custom paragraph style.\nconst bylineStyle = context.document.addStyle(\"byline\", Word.StyleType.paragraph);\nawait context.sync();\n\n// Re-acquire the style from the styles collection so property writes\n// (basedOn / quick style) land on the freshly created style.\nconst styles = context.document.getStyles();\nconst byline = styles.getByName(\"byline\");\n\n// \"Based on\" -> <w:basedOn w:val=\"Author\"/>\nbyline.baseStyle = \"Author\";\n\n// <w:qFormat/> (shown in the Quick Style gallery).\nbyline.quickStyle = true;\n\nawait context.sync();\n", "ps1": "# Add a new paragraph style \"byline\" based on the existing \"Author\"\n# style, matching the commit's styles.xml / stylesWithEffects.xml\n# addition:\n#\n#   <w:style w:type=\"paragraph\" w:customStyle=\"1\" w:styleId=\"byline\">\n#     <w:name w:val=\"byline\"/>\n#     <w:basedOn w:val=\"Author\"/>\n#     <w:qFormat/>\n#     <w:rsid w:val=\"00076B7A\"/>\n#   </w:style>\n\n$d = $word.ActiveDocument\n\n# wdStyleTypeParagraph = 1\n$byline = $d.Styles.Add(\"byline\", 1)\n\n# \"Based on\" -> <w:basedOn w:val=\"Author\"/>\n$byline.BaseStyle = \"Author\"\n\n# <w:qFormat/> (shown in the Quick Style gallery).\n$byline.QuickStyle = $true\n"}
